# Apply edit: rename headers, add "PO Forecast" sheet with data

$wb = $excel.ActiveWorkbook

# 1) Rename header on "Weekly Quantity" sheet: B1 "Requested quantity" -> "Weekly_PO_Qty"
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# 2) Rename header on "Monthly Trend" sheet: B1 "Requested quantity" -> "Monthly_PO_Qty"
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 3) Add a new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match outline/page-setup markers + page margins used by the other sheets
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1
$wsForecast.PageSetup.LeftMargin = $wsWeekly.PageSetup.LeftMargin
$wsForecast.PageSetup.RightMargin = $wsWeekly.PageSetup.RightMargin
$wsForecast.PageSetup.TopMargin = $wsWeekly.PageSetup.TopMargin
$wsForecast.PageSetup.BottomMargin = $wsWeekly.PageSetup.BottomMargin
$wsForecast.PageSetup.HeaderMargin = $wsWeekly.PageSetup.HeaderMargin
$wsForecast.PageSetup.FooterMargin = $wsWeekly.PageSetup.FooterMargin

# Header row: ds | PO_Forecast | yhat_lower | yhat_upper
$wsForecast.Range("A1").Value2 = "ds"
$wsForecast.Range("B1").Value2 = "PO_Forecast"
$wsForecast.Range("C1").Value2 = "yhat_lower"
$wsForecast.Range("D1").Value2 = "yhat_upper"

# Copy the bold/bordered header style from an existing sheet's header cell
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122) # xlPasteFormats

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @"
44934.99999999999	445	-289.9409588972224	1170.761440855442
44941.99999999999	444	-272.7049641612792	1125.04387693149
44948.99999999999	444	-267.4164165108588	1177.383370957601
44955.99999999999	444	-300.0859986122487	1141.682790943038
44969.99999999999	443	-264.0231526242285	1156.61178402286
44983.99999999999	443	-226.7208667268627	1201.7398784694
44990.99999999999	443	-285.2828192877815	1149.209291684685
44997.99999999999	443	-322.8281899918424	1165.410275825284
45004.99999999999	442	-295.3394124853311	1187.811057149449
45011.99999999999	442	-272.2445277506508	1192.59236598972
45018.99999999999	442	-259.4499077277126	1163.48216222271
45025.99999999999	442	-277.6303417628422	1139.363541969878
45039.99999999999	441	-236.1343642208084	1209.64983381655
45046.99999999999	441	-293.5230963837663	1172.343190964977
45053.99999999999	441	-235.4071189153025	1146.260087273062
45060.99999999999	441	-222.1954638793226	1145.693420719175
45067.99999999999	440	-287.2940438139869	1209.069025946299
45081.99999999999	440	-287.2891858040741	1142.290619836013
45088.99999999999	440	-287.7944019080864	1150.05362590087
45095.99999999999	440	-255.3802885461564	1124.330107436275
45102.99999999999	439	-277.5542107735411	1155.784923018085
45109.99999999999	439	-267.5269972890526	1114.73655226225
45116.99999999999	439	-324.6194692882771	1182.454140646808
45123.99999999999	439	-240.9428942336119	1161.477872325344
45130.99999999999	438	-283.189072081649	1184.90817639937
45137.99999999999	438	-315.9947911123787	1095.072181497601
45151.99999999999	438	-268.0415365312645	1136.156492554238
45165.99999999999	437	-322.1117810751449	1198.599044169901
45172.99999999999	437	-292.6633915263756	1165.747446948848
45179.99999999999	437	-260.5856831749215	1170.376999604603
45186.99999999999	437	-267.3305746432385	1101.692841329283
45193.99999999999	437	-270.3589020215453	1150.260294127405
45200.99999999999	436	-277.559285699939	1222.976595590683
45214.99999999999	436	-255.0137589501467	1185.613732181555
45221.99999999999	436	-301.4605728614911	1174.54940367599
45228.99999999999	435	-303.733369700765	1157.925395357676
45235.99999999999	435	-265.4531056192395	1174.96711300128
45249.99999999999	435	-265.4477751584012	1161.377419664424
45256.99999999999	435	-233.8646644867202	1140.001547907186
45263.99999999999	434	-292.25500030979	1148.233255855973
45298.99999999999	433	-300.1546114286565	1210.645039265828
45305.99999999999	433	-274.3922962494012	1131.563875615925
45312.99999999999	433	-278.7485430109724	1135.364406132843
45333.99999999999	432	-257.6430405722757	1150.538606578244
45347.99999999999	432	-304.3963511850745	1150.71091804669
45354.99999999999	432	-268.6026753639314	1131.248647130403
45361.99999999999	431	-307.5047548940352	1138.495372553598
45389.99999999999	430	-221.3478354984539	1124.870803562139
45396.99999999999	430	-279.0346520829982	1184.61770948408
45424.99999999999	429	-296.8648814419495	1150.734661345675
45438.99999999999	429	-252.7803788169764	1166.386281406363
45459.99999999999	428	-341.2218369918426	1115.316386214578
45466.99999999999	428	-308.4127081123312	1105.547589357853
45473.99999999999	428	-357.1705234293216	1137.943748141013
45487.99999999999	427	-360.6198226234575	1158.871411361449
45494.99999999999	427	-309.7722825454156	1170.669213845321
45501.99999999999	427	-259.931261529991	1239.455427029765
45515.99999999999	427	-297.4966336919891	1179.992451141247
45529.99999999999	426	-321.4190844288645	1157.36347353476
45543.99999999999	426	-291.9211550387095	1158.024698614707
45557.99999999999	425	-293.9558072217337	1133.026220681679
45564.99999999999	425	-298.5516582651163	1135.484846889625
45571.99999999999	425	-353.5755395998383	1132.763121665733
45578.99999999999	425	-361.2469930843484	1135.837922402386
45585.99999999999	424	-281.3460587909236	1188.599110399139
45613.99999999999	424	-341.3143876699858	1135.529793030595
45620.99999999999	423	-316.1159139548205	1139.89194758739
45627.99999999999	423	-313.1568609574044	1174.207007610303
45634.99999999999	423	-282.6540446506584	1122.784308938292
45641.99999999999	423	-259.2966762652333	1128.491520893191
45648.99999999999	422	-236.1988337576491	1146.735805268781
45655.99999999999	422	-330.2596461055865	1137.033407572068
45662.99999999999	422	-331.2953536989633	1118.547467434849
45669.99999999999	422	-301.6729839802035	1086.944740013963
"@
$lines = $data -split "`n"
$n = $lines.Count
$arr = New-Object 'object[,]' $n,4
for ($i = 0; $i -lt $n; $i++) {
    $parts = $lines[$i].Trim() -split "`t"
    $arr[$i,0] = [double]$parts[0]
    $arr[$i,1] = [double]$parts[1]
    $arr[$i,2] = [double]$parts[2]
    $arr[$i,3] = [double]$parts[3]
}

$dataRange = $wsForecast.Range("A2").Resize($n, 4)
$dataRange.Value = $arr

# Copy the date number-format style used for column A on the other sheets
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2").Resize($n, 1).PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0
$null = $wsForecast.Range("A1").Select()

Write-Output "Edit applied: renamed headers + added PO Forecast sheet with $n data rows"
